$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.940.27"
$ws.Range("E2").Value = "  +2.01%  "

$ws.Range("D3").Value = "3.468.43"
$ws.Range("E3").Value = "  +2.07%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'577.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").Value = "'147.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.29%  "

$ws.Range("D7").Value = "3.469.15"
$ws.Range("E7").Value = "  +2.06%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +1.19%  "

$ws.Range("D10").Value = "'7.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.27%  "

$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("D12").Value = "'0.402"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.16%  "

$ws.Range("D13").Value = "4.063.31"
$ws.Range("E13").Value = "  +2.13%  "

$ws.Range("E14").Value = "  +6.12%  "

$ws.Range("E15").Value = "  +2.69%  "

$ws.Range("D16").Value = "3.475.51"
$ws.Range("E16").Value = "  +2.38%  "

$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").Value = "62.925.70"
$ws.Range("E18").Value = "  +1.98%  "

$ws.Range("E19").Value = "  +2.83%  "

$ws.Range("D20").Value = "'14.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.07%  "

$ws.Range("D21").Value = "'9.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("D22").Value = "'387.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'74.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "'0.556"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").Value = "3.608.86"
$ws.Range("E26").Value = "  +2.10%  "

$ws.Range("D27").Value = "'0.0000114"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.85%  "

$ws.Range("E28").Value = "  -0.82%  "

$ws.Range("D29").Value = "'7.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.67%  "

$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("D31").Value = "'8.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.66%  "

$ws.Range("D32").Value = "'2.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.98%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'23.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.09%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("E36").Value = "  +3.38%  "

$ws.Range("D37").Value = "'31.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +17.70%  "

$ws.Range("E38").Value = "  +1.37%  "

$ws.Range("D39").Value = "'169.76"
$ws.Range("D39").Style = "Normal"

$ws.Range("E40").Value = "  +5.40%  "

$ws.Range("D41").Value = "3.506.72"
$ws.Range("E41").Value = "  +2.22%  "

$ws.Range("E42").Value = "  -1.19%  "

$ws.Range("D43").Value = "'0.798"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.90%  "

$ws.Range("D44").Value = "'42.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("E46").Value = "  +2.16%  "

$ws.Range("E47").Value = "  +3.48%  "

$ws.Range("D48").Value = "2.614.32"
$ws.Range("E48").Value = "  +5.42%  "

$ws.Range("D49").Value = "'2.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.86%  "

$ws.Range("D50").Value = "'22.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "

$ws.Range("E51").Value = "  +0.90%  "
